# Update the "想去人数" (want-to-go count) figures in column F that changed
# between crawls, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 11499
    "F3"  = 10963
    "F5"  = 7
    "F11" = 10623
    "F12" = 4107
    "F13" = 9
    "F14" = 2455
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
